$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("B","C","D","E","G","H","K","L")

$bVals = @(15.81101989289832, 15.70406812970731, 15.64355627489336, 15.62022006136597, 15.61642571043213, 15.64323616581842, 15.77308903402618, 16.0673801499365, 16.30591100603259, 16.41883524718436, 16.46219293071936, 16.45282921879013, 16.42239059922462, 16.40382246395017, 16.2986166656455, 16.23517968965044, 16.1991104442778, 16.18697098855047, 16.24188968034645, 16.43131530981455, 16.55856933152689, 16.49034870329508, 16.23885483991041, 15.98370545594099)
$cVals = @(6.783469263426012, 6.727077132321249, 6.690948633323608, 6.675848203110776, 6.67331798088866, 6.690746512791337, 6.764337465716776, 6.896711608306524, 6.986686478668654, 7.026037929165688, 7.040712253289068, 7.037561983587526, 7.027249749500173, 7.020903612313119, 6.984082984236685, 6.961090242431703, 6.947716853526598, 6.943163343862367, 6.963553237343692, 7.030284870748618, 7.072573389245321, 7.050124461080473, 6.962440198763474, 6.862179340597359)
$dVals = @(6.044966231151203, 5.936878113912289, 5.871396535223088, 5.844973347909519, 5.840602690138511, 5.871039074900001, 6.007536929699472, 6.280539739588145, 6.482133361136261, 6.573609758408057, 6.608180339508592, 6.600738603766082, 6.576455498203858, 6.561571281338073, 6.476147384676654, 6.423657302423165, 6.393446879805803, 6.383215937162849, 6.429247250797711, 6.583590205491316, 6.684041119090535, 6.630478752278816, 6.426720135699177, 6.206347970868531)
$eVals = @(11.56995848728289, 11.56173579769607, 11.55835612323435, 11.55739948661268, 11.557266064147, 11.55834151761518, 11.56677747765337, 11.59651803813078, 11.62634357310261, 11.64162474095863, 11.64765571210227, 11.64634600305599, 11.6421160284537, 11.63955680065029, 11.6253792191591, 11.61711917018286, 11.6125296161391, 11.61100345331866, 11.61798177569076, 11.64335186102497, 11.66135569364427, 11.65161723467602, 11.61759129553309, 11.58706600878184)
$gVals = @(63.3751526831309, 62.56557699248685, 62.07145728175877, 61.87102197081327, 61.83780069556079, 62.06875017338746, 63.09550300682747, 65.1246185259731, 66.61443576404328, 67.28977238656225, 67.54497974995351, 67.49004265128272, 67.31078007269517, 67.20090240775622, 66.57023699967688, 66.18259989510325, 65.9594272579223, 65.88383376412656, 66.22388802694755, 67.36344961533585, 68.10505920907652, 67.70959715017601, 66.20522262928755, 64.5751018599542)
$hVals = @(22.40348432731636, 22.30166302174474, 22.24131106177493, 22.21727239562632, 22.2133146189798, 22.24098460663279, 22.36793307762791, 22.63358729583828, 22.83826612496435, 22.93328532621208, 22.96952688597961, 22.96171030104787, 22.93626181969961, 22.92070728350394, 22.8320935747188, 22.77821072192606, 22.74739940443859, 22.73699869464753, 22.78392805009755, 22.94372972341207, 23.04967677924213, 22.9929979150697, 22.78134272574688, 22.56000002298026)
$kVals = @(12.12358978834094, 12.073333964015, 12.0467015912705, 12.03692153788667, 12.03536263810614, 12.04656533755967, 12.10539149124784, 12.25373251606875, 12.3819933347429, 12.44431952226718, 12.46847362548043, 12.46324735994394, 12.44629573858441, 12.43598373073312, 12.37799859237003, 12.34343221315477, 12.32392639549916, 12.31738719549668, 12.34707310365268, 12.4512600083071, 12.52256278957121, 12.48422019649977, 12.34542591286192, 12.2101512194327)
$lVals = @(10.02090689288339, 10.01898293633741, 10.01975977099589, 10.0205692754815, 10.02073347657756, 10.01976869182987, 10.01983749279802, 10.03547675628771, 10.05636197004015, 10.06788604482039, 10.07253894375041, 10.07152404016533, 10.06826306214797, 10.0663031902635, 10.05564939630183, 10.0496305350051, 10.04635922042535, 10.04528439662328, 10.05025154268747, 10.06921306447589, 10.08328883171742, 10.07562302111448, 10.04997019623729, 10.02959093542027)

$data = @(
    $bVals,
    $cVals,
    $dVals,
    $eVals,
    $gVals,
    $hVals,
    $kVals,
    $lVals
)

for ($ci = 0; $ci -lt $colLetters.Length; $ci++) {
    $colVals = $data[$ci]
    $colLetter = $colLetters[$ci]
    for ($ri = 0; $ri -lt $colVals.Length; $ri++) {
        $rowNum = $ri + 2
        $ws.Range("$colLetter$rowNum").Value2 = $colVals[$ri]
    }
}

Write-Output "applied 192 cell updates"
